$wb = $excel.ActiveWorkbook

# Mapping of row -> new "想去人数" (F column) value.
# Same update is applied to both the "展览" sheet and the "全部类型" sheet,
# which mirror each other's data.
$updates = @{
    2  = 8428
    3  = 8030
    4  = 132
    5  = 195
    8  = 136
    9  = 136
    13 = 144
    14 = 2085
    20 = 60
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
